$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")
$rapport = $wb.Worksheets.Item("Rapport")

# ---------------------------------------------------------------------------
# 0) Rows 59-71 already carry the same per-column formatting as the rest of
#    the table for columns A-E and G-L (column F is left with the default
#    column formatting until a Type is actually entered). Re-applying that
#    formatting (copied from row 59, which already has it) is a no-op for
#    rows 59-61 and brings rows 62-71 up to the same state, so the later
#    value assignments preserve the correct number formats instead of
#    falling back to the raw column defaults.
# ---------------------------------------------------------------------------
$ws.Range("A59:E59").Copy()
$ws.Range("A59:E71").PasteSpecial(-4122) | Out-Null
$ws.Range("G59:L59").Copy()
$ws.Range("G59:L71").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 1) New rows of the "Journal" sheet for the work done on 11.05.2021
#    (date serial 44327), continuing the existing table that ran out at
#    row 58. Columns: A=Date, B=Semaine, C=Debut, D=Fin, E=Temps (formula),
#    F=Type, G=Titre, H=Description.
# ---------------------------------------------------------------------------

function Set-JournalRow($Row, $Start, $End, $Type, $Title, $Description) {
    $ws.Cells.Item($Row, 1).Value = 44327
    $ws.Cells.Item($Row, 2).Value = 2
    $ws.Cells.Item($Row, 3).Value = $Start
    $ws.Cells.Item($Row, 4).Value = $End
    $ws.Cells.Item($Row, 6).Value = $Type
    $ws.Cells.Item($Row, 7).Value = $Title
    if ($Description) {
        $ws.Cells.Item($Row, 8).Value = $Description
    }
}

# Row 59 - 08:00 -> 08:44
Set-JournalRow 59 0.33333333333333331 0.36388888888888887 "Réalisation" "Environnement" "Installation du serveur MySQL."

# Row 60 - 08:44 -> 09:35
Set-JournalRow 60 0.36388888888888887 0.39930555555555558 "Communication" "Sprint Review" $null

# Rows 62-64 (Review Scénarios block) were filled in before row 61 was
# revisited, matching the order the new shared strings were recorded in
# the workbook.

# Row 62 - 10:35 -> 12:15
Set-JournalRow 62 0.44097222222222227 0.51041666666666663 "Analyse" "Review Scénarios" "Review et corrections des scénarios."

# Row 63 - 13:30 -> 14:15
Set-JournalRow 63 0.5625 0.59375 "Analyse" "Review Scénarios" "Review et corrections des scénarios."

# Row 64 - 14:20 -> 14:34
Set-JournalRow 64 0.59722222222222221 0.6069444444444444 "Analyse" "Review Scénarios" "Review et corrections des scénarios."

# Row 61 (added afterwards) - 09:55 -> 10:35
Set-JournalRow 61 0.41319444444444442 0.44097222222222227 "Conception" "Review MCD MLD" "Review et corrections des diagrammes MCD et MLD."

# Row 65 - 14:34 -> 14:54
Set-JournalRow 65 0.6069444444444444 0.62083333333333335 "Analyse" "Corrections" "Corrections de certains détails."

# Row 66 - 14:54 -> 15:05
Set-JournalRow 66 0.62083333333333335 0.62847222222222221 "Réalisation" "Environnement" "Mise en place de la base de données."

# Row 67 - 15:20 -> 15:49
Set-JournalRow 67 0.63888888888888895 0.65902777777777777 "Réalisation" "Environnement" "Mise en place de la base de données."

# Row 68 - 15:49 -> 16:05
Set-JournalRow 68 0.65902777777777777 0.67013888888888884 "Réalisation" "Environnement" "Création du script de création de la base de données."

# Fill down the "Temps" (E) formula over the new rows, extending the
# existing shared formula range from E25:E58 up to E25:E68.
$ws.Range("E25:E68").Formula = "=D25-C25"

# Rows 61 and 68 wrap their description text onto two lines, same as the
# other multi-line rows in the sheet.
$ws.Rows(61).RowHeight = 30
$ws.Rows(68).RowHeight = 30

# ---------------------------------------------------------------------------
# 2) Selection bookkeeping left over from the editing session.
# ---------------------------------------------------------------------------
$ws.Range("H38").Select()

$rapport.Range("C9").Select()
